# issue #5: stock data output to json file
# Adds a "property_category" column (value "stock") to the 股票 (stock) sheet,
# and fixes a handful of company-name strings that had a stray internal space.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Insert a new column before the existing "date" column (H), shifting
# date/legislator_name/legislator_id one column to the right (H->I, I->J, J->K).
$ws.Columns.Item(8).Insert()

# New header + values for the inserted property_category column.
$ws.Range("H1").Value = "property_category"
$ws.Range("H2:H15").Value = "stock"

# Fix company names that had an extraneous space inserted mid-word.
$ws.Range("B5").Value = "鴻海精密工業股份有限公司"
$ws.Range("B6").Value = "台灣積體電路製造股份有限公司"
$ws.Range("B9").Value = "彰化商業銀行股份有限公司"
$ws.Range("B10").Value = "中華票券金融股份有限公司"
$ws.Range("B11").Value = "安泰商業銀行股份有限公司"
$ws.Range("B13").Value = "台灣塑膠工業股份有限公司"
$ws.Range("B14").Value = "裕隆汽車製造股份有限公司"
